$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": clear B and C detail columns to reflect re-run with empty mapping ---
$ws2 = $wb.Worksheets.Item("Range Status")
for ($r = 2; $r -le 7; $r++) {
    $ws2.Cells.Item($r, 2).Value = 0
    $ws2.Cells.Item($r, 3).ClearContents()
}

# --- Sheet "Species qualification": Range Analysis row's species count now 0 ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Cells.Item(5, 2).Value = 0

# --- Sheet "High Priority break-up": only the IUCN break-up row remains (Range row removed) ---
$ws5 = $wb.Worksheets.Item("High Priority break-up")
$ws5.Cells.Item(2, 1).Value = "IUCN"
$ws5.Cells.Item(2, 2).Value = 33
$ws5.Cells.Item(2, 3).Value = 100
$ws5.Cells.Item(2, 4).Value = 33
$ws5.Cells.Item(2, 5).Value = 100

$ws5.Range("A3:E3").Delete() | Out-Null

$wb.Save()
